# Actualización 10 de Mayo
$wb = $excel.ActiveWorkbook

# --- Estadisticos 1P ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 6
$ws1.Range("H3").Value = 6.4

$ws1.Range("D5").Value = 0
$ws1.Range("E5").Value = 7
$ws1.Range("H5").Value = 6.5

$ws1.Range("D6").Value = 1
$ws1.Range("E6").Value = 9
$ws1.Range("H6").Value = 6.2

# --- Estadisticos 2P ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("D2").Value = 6
$ws2.Range("E2").Value = 8
$ws2.Range("F2").Value = 28
$ws2.Range("G2").Value = 77.78
$ws2.Range("H2").Value = 6.7

$ws2.Range("E3").Value = 10

$ws2.Range("E5").Value = 8

$ws2.Range("D6").Value = 8
$ws2.Range("E6").Value = 10
$ws2.Range("F6").Value = 18
$ws2.Range("G6").Value = 62.07
$ws2.Range("H6").Value = 6.6

# --- Estadisticos Final ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("E2").Value = 8
$ws3.Range("F2").Value = 28
$ws3.Range("G2").Value = 77.78
$ws3.Range("H2").Value = 6.5

$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 9
$ws3.Range("H3").Value = 6.1

$ws3.Range("D5").Value = 0
$ws3.Range("E5").Value = 8
$ws3.Range("H5").Value = 6.5

$ws3.Range("D6").Value = 1
$ws3.Range("E6").Value = 10
$ws3.Range("H6").Value = 6.4

# --- Rescatables ---
$ws4 = $wb.Worksheets.Item("Rescatables")

$ws4.Range("A2").Value = 18330051920088
$ws4.Range("B2").Value = "BERINSTAIN"
$ws4.Range("C2").Value = "SAN JUAN"
$ws4.Range("D2").Value = "LUIS FERNANDO"
$ws4.Range("E2").Value = "TEMAS DE ADMINISTRACIÓN"
$ws4.Range("F2").Value = "6ARHV"
$ws4.Range("G2").Value = 1
